$wb = $excel.ActiveWorkbook

# ================= 1. Update "总计" (Total) summary sheet =================
$total = $wb.Worksheets.Item("总计")

# Insert a new row at position 2, shifting existing data rows down
$total.Range("A2:D2").Insert(-4121)

# New row 2: summary for 2022-Q3
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0

# Re-number the index column for the rows that shifted down
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

# Match formatting of new row 2 to the other data rows (copy from row 4)
$total.Range("A4:D4").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

# ================= 2. Insert new "2022-Q3" detail sheet =================
$newSheet = $wb.Worksheets.Add($null, $total)
$newSheet.Name = "2022-Q3"

$q1 = $wb.Worksheets.Item("2022-Q1")

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'501307"
$newSheet.Range("C2").Value = "银河中证沪港深高股息指数（LOF）A"
$newSheet.Range("D2").Value = "'0.15"
$newSheet.Range("E2").Value = "'90.33"
$newSheet.Range("F2").Value = "'1.25"
$newSheet.Range("G2").Value = "'0.0019"
$newSheet.Range("H2").Value = 10

# Row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'501308"
$newSheet.Range("C3").Value = "银河中证沪港深高股息指数（LOF）C"
$newSheet.Range("D3").Value = "'0.01"
$newSheet.Range("E3").Value = "'90.33"
$newSheet.Range("F3").Value = "'1.25"
$newSheet.Range("G3").Value = "'0.0001"
$newSheet.Range("H3").Value = 10

# Apply the same cell formatting template used by the sibling quarter sheets
$q1.Range("A1:H3").Copy()
$newSheet.Range("A1:H3").PasteSpecial(-4122)

Write-Output "done"
